$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2024-11-15 Friday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-11-16 Saturday", 2)

# Helper to set the text of a table cell without disturbing the
# end-of-cell marker (trim the last character off the cell range
# before assigning new text).
function Set-CellText($table, $row, $col, $text) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    $r.End = $r.End - 1
    $r.Text = $text
}

$t = $d.Tables(1)

# Row 1
Set-CellText $t 1 1 "61÷6="
Set-CellText $t 1 2 "81÷6="
Set-CellText $t 1 3 "10÷7="
Set-CellText $t 1 4 "60÷2="
Set-CellText $t 1 5 "57÷4="

# Row 5
Set-CellText $t 5 1 "40÷6="
Set-CellText $t 5 2 "75÷8="
Set-CellText $t 5 3 "61÷5="
Set-CellText $t 5 4 "54÷9="
Set-CellText $t 5 5 "72÷5="

# Row 9
Set-CellText $t 9 1 "56÷7="
Set-CellText $t 9 2 "83÷5="
Set-CellText $t 9 3 "29÷2="
Set-CellText $t 9 4 "86÷7="
Set-CellText $t 9 5 "45÷2="

# Row 13
Set-CellText $t 13 1 "24÷4="
Set-CellText $t 13 2 "57÷9="
Set-CellText $t 13 3 "33÷8="
Set-CellText $t 13 4 "66÷6="
Set-CellText $t 13 5 "85÷9="

# Row 17
Set-CellText $t 17 1 "87÷6="
Set-CellText $t 17 2 "70÷8="
Set-CellText $t 17 3 "52÷5="
Set-CellText $t 17 4 "69÷2="
Set-CellText $t 17 5 "97÷4="
